$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("N3").Value = 1.93
$ws.Range("O3").Value = 1.97
# Row 11
$ws.Range("J11").Value = 1.08
$ws.Range("K11").Value = 8
$ws.Range("L11").Value = 1.44
$ws.Range("M11").Value = 2.63
$ws.Range("N11").Value = 2.4
$ws.Range("O11").Value = 1.53
# Row 13
$ws.Range("K13").Value = 19
$ws.Range("L13").Value = 1.14
# Row 15
$ws.Range("J15").Value = 1.06
$ws.Range("K15").Value = 10
$ws.Range("L15").Value = 1.33
# Row 16
$ws.Range("H16").Value = 5.75
$ws.Range("J16").Value = 1.04
$ws.Range("L16").Value = 1.2
$ws.Range("R16").Value = 2.38
$ws.Range("S16").Value = 1.53
$ws.Range("T16").Value = 26
$ws.Range("Y16").Value = 81
$ws.Range("AG16").Value = 9.5
$ws.Range("AH16").Value = 7
# Row 17
$ws.Range("J17").Value = 1.05
$ws.Range("L17").Value = 1.3
$ws.Range("N17").Value = 2.03
$ws.Range("O17").Value = 1.83
# Row 19
$ws.Range("K19").Value = 13
# Row 21
$ws.Range("G21").Value = 2.67
$ws.Range("H21").Value = 3
$ws.Range("I21").Value = 2.62
$ws.Range("M21").Value = 2.67
$ws.Range("O21").Value = 1.62
$ws.Range("R21").Value = 1.75
$ws.Range("S21").Value = 1.85
$ws.Range("U21").Value = 13
$ws.Range("Z21").Value = 7.9
$ws.Range("AA21").Value = 5.8
$ws.Range("AE21").Value = 7.6
$ws.Range("AF21").Value = 13
$ws.Range("AG21").Value = 9.75
# Row 28
$ws.Range("H28").Value = 4.05
$ws.Range("I28").Value = 1.95
$ws.Range("L28").Value = 1.08
$ws.Range("M28").Value = 6.6
$ws.Range("N28").Value = 1.27
$ws.Range("O28").Value = 3.45
$ws.Range("P28").Value = 1.18
$ws.Range("Q28").Value = 4.35
$ws.Range("R28").Value = 1.28
$ws.Range("S28").Value = 3.35
$ws.Range("U28").Value = 32
$ws.Range("Y28").Value = 18
$ws.Range("Z28").Value = 35
$ws.Range("AA28").Value = 10.5
$ws.Range("AC28").Value = 22
$ws.Range("AD28").Value = 75
$ws.Range("AE28").Value = 18
$ws.Range("AF28").Value = 17
$ws.Range("AG28").Value = 9.75
$ws.Range("AH28").Value = 23
$ws.Range("AJ28").Value = 14.5
# Row 29
$ws.Range("G29").Value = 1.91
$ws.Range("H29").Value = 3.7
$ws.Range("I29").Value = 3.5
$ws.Range("K29").Value = 9.5
$ws.Range("M29").Value = 4.9
$ws.Range("O29").Value = 2.6
$ws.Range("P29").Value = 1.26
$ws.Range("Q29").Value = 3.5
$ws.Range("R29").Value = 1.4
$ws.Range("S29").Value = 2.7
$ws.Range("U29").Value = 13.5
$ws.Range("V29").Value = 8.5
$ws.Range("W29").Value = 20
$ws.Range("Y29").Value = 16
$ws.Range("Z29").Value = 9.5
$ws.Range("AA29").Value = 8
$ws.Range("AB29").Value = 10.75
$ws.Range("AC29").Value = 30
$ws.Range("AE29").Value = 16.5
$ws.Range("AI29").Value = 26
# Row 34
$ws.Range("G34").Value = 1.91
# Row 35
$ws.Range("H35").Value = 4.5
$ws.Range("R35").Value = 2.12
# Row 40
$ws.Range("K40").Value = 13
# Row 41
$ws.Range("K41").Value = 13
# Row 42
$ws.Range("G42").Value = 1.75
$ws.Range("I42").Value = 4.2
$ws.Range("K42").Value = 13
$ws.Range("U42").Value = 9
$ws.Range("Z42").Value = 13
# Row 44
$ws.Range("G44").Value = 2.05
$ws.Range("I44").Value = 3.9
$ws.Range("T44").Value = 5.5
# Row 46 (previously blank odds, now populated)
$ws.Range("G46").Value = 1.7
$ws.Range("H46").Value = 3.5
$ws.Range("I46").Value = 4.4
$ws.Range("J46").Value = 1.07
$ws.Range("K46").Value = 7.1
$ws.Range("L46").Value = 1.32
$ws.Range("M46").Value = 3.1
$ws.Range("N46").Value = 1.95
$ws.Range("O46").Value = 1.75
$ws.Range("P46").Value = 1.44
$ws.Range("Q46").Value = 2.62
$ws.Range("R46").Value = 1.91
$ws.Range("S46").Value = 1.8
$ws.Range("T46").Value = 6.3
$ws.Range("U46").Value = 7.6
$ws.Range("V46").Value = 8.25
$ws.Range("W46").Value = 13
$ws.Range("X46").Value = 14.5
$ws.Range("Y46").Value = 29
$ws.Range("Z46").Value = 7.1
$ws.Range("AA46").Value = 6.9
$ws.Range("AB46").Value = 16.5
$ws.Range("AC46").Value = 90
$ws.Range("AD46").Value = 700
$ws.Range("AE46").Value = 11.75
$ws.Range("AF46").Value = 25
$ws.Range("AG46").Value = 15
$ws.Range("AH46").Value = 75
$ws.Range("AI46").Value = 45
$ws.Range("AJ46").Value = 50
# Row 47
$ws.Range("G47").Value = 1.52
$ws.Range("H47").Value = 3.9
$ws.Range("I47").Value = 5.4
$ws.Range("K47").Value = 7.4
$ws.Range("L47").Value = 1.3
$ws.Range("M47").Value = 3.25
$ws.Range("N47").Value = 1.88
$ws.Range("O47").Value = 1.82
$ws.Range("P47").Value = 1.4
$ws.Range("Q47").Value = 2.72
$ws.Range("R47").Value = 2
$ws.Range("S47").Value = 1.72
$ws.Range("T47").Value = 6.2
$ws.Range("U47").Value = 6.7
$ws.Range("W47").Value = 10.25
$ws.Range("X47").Value = 13
$ws.Range("Z47").Value = 7.4
$ws.Range("AA47").Value = 7.8
$ws.Range("AB47").Value = 19.5
$ws.Range("AE47").Value = 13.5
$ws.Range("AF47").Value = 32
$ws.Range("AG47").Value = 18
$ws.Range("AH47").Value = 110
$ws.Range("AI47").Value = 60
$ws.Range("AJ47").Value = 65
# Row 51
$ws.Range("N51").Value = 2.2
$ws.Range("O51").Value = 1.65
# Row 52
$ws.Range("J52").Value = 26
$ws.Range("K52").Value = 1.02
$ws.Range("L52").Value = 1.05
# Row 54
$ws.Range("J54").Value = 1.03
$ws.Range("K54").Value = 8
$ws.Range("L54").Value = 1.38
$ws.Range("M54").Value = 2.62
# Row 55
$ws.Range("J55").Value = 1.05
$ws.Range("L55").Value = 1.38
$ws.Range("M55").Value = 2.62
